$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.688.15'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.632.15'
$ws.Range("E3").Value = '  -1.03%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.95'
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("E6").Value = '  -1.64%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.249'
$ws.Range("E8").Value = '  -1.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0620'
$ws.Range("E9").Value = '  -1.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.99'
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.859.44'
$ws.Range("E12").Value = '  -0.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.634.62'
$ws.Range("E13").Value = '  -1.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.11'
$ws.Range("E14").Value = '  -2.22%  '
$ws.Range("E15").Value = '  -2.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.96'
$ws.Range("E16").Value = '  -2.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.672.50'
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0721'
$ws.Range("E18").Value = '  -3.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '211.12'
$ws.Range("E19").Value = '  -3.32%  '
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.30'
$ws.Range("E21").Value = '  -1.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.34'
$ws.Range("E22").Value = '  -7.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.17'
$ws.Range("E23").Value = '  -2.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.15'
$ws.Range("E24").Value = '  -3.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.52'
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -2.65%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.00'
$ws.Range("E28").Value = '  -2.82%  '
$ws.Range("E29").Value = '  -2.10%  '
$ws.Range("E30").Value = '  -3.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.19'
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.37'
$ws.Range("E32").Value = '  +0.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.94'
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.259.41'
$ws.Range("E34").Value = '  -1.91%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  -2.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0173'
$ws.Range("E37").Value = '  -3.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.523'
$ws.Range("E38").Value = '  -3.29%  '
$ws.Range("E39").Value = '  -0.11%  '
$ws.Range("E40").Value = '  -3.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.797'
$ws.Range("E41").Value = '  -2.50%  '
$ws.Range("B42").Value = 'RocketPoolETH'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.771.12'
$ws.Range("E42").Value = '  -1.71%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.15'
$ws.Range("E43").Value = '  -4.58%  '
$ws.Range("E44").Value = '  -3.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.90'
$ws.Range("E45").Value = '  -1.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.79'
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.405'
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0953'
$ws.Range("E51").Value = '  -2.81%  '
